$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows whose x_nrSteps (F) value is -2 / alienID (H) value is 36,
# decrementing x_corrSteps (D) and x_nrSteps (F) by 1, and bumping alienID (H) by 10.
$rows = @(2, 9, 12, 20, 22, 28)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = $ws.Range("D$r").Value() - 1
    $ws.Range("F$r").Value = $ws.Range("F$r").Value() - 1
    $ws.Range("H$r").Value = $ws.Range("H$r").Value() + 10
}

# Update the selected cell shown in the saved sheet view.
$ws.Range("E28").Select()
